$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert two new columns before column C (shifts old C,D,E... to the right by two)
$ws.Range("C1:D1").EntireColumn.Insert()

# 2. Set the width of the two newly inserted columns (Phones / E-mail)
$ws.Range("C1:D1").ColumnWidth = 19.7109375

# 3. Update existing labels whose text changed
$ws.Range("A6").Value = "Группировка: {{GroupingTitle}}"
$ws.Range("A8").Value = "{{GroupingTitle}}"

# 4. New column headers for the inserted columns
$ws.Range("C8").Value = "Телефоны"
$ws.Range("D8").Value = "E-mail"

# 5. New data-row placeholders for the inserted columns
$ws.Range("C11").Value = "{{item.Phones}}"
$ws.Range("D11").Value = "{{item.Emails}}"

# 6. Remove the old "Warehouse residue" column content (now pushed out to column I)
$ws.Range("I8").ClearContents()
$ws.Range("I9").ClearContents()
$ws.Range("I11").ClearContents()

# 7. Remove the old "№" / "Номенклатура" header labels (now blank in the new layout)
$ws.Range("A10").ClearContents()
$ws.Range("B10").ClearContents()

# 8. Merge the "Группировка" row across the widened header block
$ws.Range("A6:G6").Merge()

# 9. Restore selection to match the new layout
$ws.Range("A8:B8").Select()
